$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-coerced to a number by Excel
# (e.g. "1.000", "0.9996") need NumberFormat forced to Text ("@") beforehand so the
# literal string (including insignificant trailing/leading zeros) is preserved.

$ws.Range("D2").Value = "27.163.60"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.867.47"
$ws.Range("E3").Value = "  -2.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.04"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5176"
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3754"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07171"
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.70"
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8849"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07565"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "1.857.33"
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.337"
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.37"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008562"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.16"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "27.215.47"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.033"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "2.118.10"
$ws.Range("E22").Value = "  -3.29%  "
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.476"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.79"
$ws.Range("E25").Value = "  -2.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.847"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.01"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.137"
$ws.Range("E28").Value = "  -4.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.66"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.753"
$ws.Range("E30").Value = "  -3.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.694"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09007"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05163"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.100"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7522"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.174"
$ws.Range("E36").Value = "  -4.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02036"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.541"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.026"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  -1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5334"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.642"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.84"
$ws.Range("E43").Value = "  +3.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.485"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1484"
$ws.Range("E45").Value = "  -1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4660"
$ws.Range("E46").Value = "  -3.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9994"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.14"
$ws.Range("E48").Value = "  -3.69%  "
$ws.Range("E49").Value = "  -3.64%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "65.02"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.44"
$ws.Range("E51").Value = "  -1.58%  "
